{"js": "// The edit removes 4 consecutive paragraphs that followed the\n// \"LOT2013: Engenharia Bioqu\u00edmica I (Requisito fraco)\" requirement line:\n//   1) an empty paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) an empty paragraph\n//   4) an empty paragraph with a page-break-before\n// (A later empty paragraph and a final page-break-before paragraph stay.)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (the last requirement line) by its text.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"LOT2013: Engenharia Bioqu\u00edmica I (Requisito fraco)\") {\n    anchorIndex = i;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOT2013' requirement paragraph.\");\n}\n\n// The four paragraphs that must be removed sit right after the anchor.\nconst toRemoveTexts = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\",\n  \"\"\n];\n\nconst removeStart = anchorIndex + 1;\nfor (let k = 0; k < toRemoveTexts.length; k++) {\n  const idx = removeStart + k;\n  if (idx >= items.length || items[idx].text !== toRemoveTexts[k]) {\n    throw new Error(\"Unexpected document shape near index \" + idx + \"; refusing to delete.\");\n  }\n}\n\n// Delete from the last one back to the first so indices stay valid.\nfor (let k = toRemoveTexts.length - 1; k >= 0; k--) {\n  items[removeStart + k].delete();\n}\n\nawait context.sync();\n", "ps1": "# The edit removes 4 consecutive paragraphs that followed the\n# \"LOT2013: Engenharia Bioqu\u00edmica I (Requisito fraco)\" requirement line:\n#   1) an empty paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) an empty paragraph\n#   4) an empty paragraph with a page-break-before\n# (A later empty paragraph and a final page-break-before paragraph stay.)\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (the last requirement line) by its text.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -eq \"LOT2013: Engenharia Bioqu\u00edmica I (Requisito fraco)`r\") {\n        $anchorIndex = $i\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the 'LOT2013' requirement paragraph.\"\n}\n\n# The four paragraphs that must be removed sit right after the anchor.\n$toRemoveTexts = @(\n    \"`r\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx`r\",\n    \"`r\",\n    \"`r\"\n)\n\n$removeStart = $anchorIndex + 1\nfor ($k = 0; $k -lt $toRemoveTexts.Length; $k++) {\n    $idx = $removeStart + $k\n    if ($idx -gt $d.Paragraphs.Count -or $d.Paragraphs.Item($idx).Range.Text -ne $toRemoveTexts[$k]) {\n        throw \"Unexpected document shape near paragraph $idx; refusing to delete.\"\n    }\n}\n\n# Delete from the last one back to the first so indices stay valid.\nfor ($k = $toRemoveTexts.Length - 1; $k -ge 0; $k--) {\n    $d.Paragraphs.Item($removeStart + $k).Range.Delete()\n}\n"}
